$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update product names (A2, A3)
$ws.Range("A2").Value = "Coba 3"
$ws.Range("A3").Value = "Coba 4"

# Clear the "Kode Produk" (B) and "Nomor Seri" (C) values for both rows
$ws.Range("B2:C2").ClearContents()
$ws.Range("B3:C3").ClearContents()

# Update the I1 header value back to "Agen"
$ws.Range("I1").Value = "Agen"

# Update J3 "Keterangan" value back to "Tes"
$ws.Range("J3").Value = "Tes"

# Update the selection to C6 (also resets the view's top-left cell to default)
[void]$ws.Range("C6").Select()
